$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'AU-8 b,AU-12 (3),AU-7 b,AU-7 a,AC-6 (9),AC-6 (8),CM-5 (1)'
$ws.Range("A8").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A10").Value = 'CM-7 (2),CM-7 (5) (b)'
$ws.Range("A15").Value = 'IA-2,AU-3 (1),IA-8'
$ws.Range("A16").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A17").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A19").Value = 'IA-5 (1) (b),CM-6 b,IA-5 (1) (a)'
$ws.Range("A21").Value = 'SC-10,AC-12,MA-4 e,MA-4 (7)'
$ws.Range("A22").Value = 'MA-4 (1) (a),AU-7 (1),AU-7 a,AU-12 a,AU-3,AU-3 (1),AU-6 (4),AU-14 (1),CM-6 b,CM-5 (1)'
$ws.Range("A25").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A29").Value = 'SC-8 (1),SC-8,SC-8 (2)'
$ws.Range("A31").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AC-2 (4),AU-3 (1),AU-12 c'
$ws.Range("A45").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a,AC-8 b'
$ws.Range("A55").Value = 'SC-8,AC-17 (2)'
$ws.Range("A56").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A65").Value = 'CM-6 b,IA-2 (2)'
$ws.Range("A67").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A69").Value = 'AU-8 b,AU-12 (3),AU-7 b,AU-7 a,AU-12 a,CM-6 b,AU-12 c,CM-5 (1)'
$ws.Range("A77").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AC-2 (4),AU-3 (1),AU-12 c'
$ws.Range("A80").Value = 'IA-2 (3),IA-2 (1),IA-2 (4),IA-2 (2)'
$ws.Range("A81").Value = 'CM-5 (3),CM-6 b'
$ws.Range("A86").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A88").Value = 'AC-2 (4),CM-5 (1),AC-6 (9),AU-12 c'
$ws.Range("A89").Value = 'IA-2 (5),IA-2 (4),IA-2 (2),IA-2 (3),IA-2'
$ws.Range("A90").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A97").Value = 'AU-8 b,AU-8 (1) (a),AU-8 (1) (b)'
$ws.Range("A101").Value = 'AC-3 (4),IA-11'
$ws.Range("A102").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A119").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A123").Value = 'CM-7 a,CM-7 b'
$ws.Range("A124").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A125").Value = 'CM-7 a,AC-18 (1)'
$ws.Range("A128").Value = 'CM-7 a,IA-5 (1) (c),CM-6 b'
$ws.Range("A136").Value = 'AC-11 (1),AC-11 b'
$ws.Range("A139").Value = 'SI-6 b,SI-6 d,CM-3 (5)'
$ws.Range("A148").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-14 (1),AU-12 c'
$ws.Range("A157").Value = 'MA-4 (1) (a),AU-12 a,AU-3,AU-3 (1),AU-12 c'
$ws.Range("A159").Value = 'SC-8,AC-17 (2)'
